# Update data sources: zero out the metric columns (B:AE) for the
# existing data rows (2025-2028 => rows 2-5), leaving the "year" column
# (A) and header row (1) untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:AE5").Value = 0
